# This script applies the "Deploying to gh-pages" content update to the
# StructureDefinition-total-dependent-count workbook:
#  - Metadata sheet: URL, Version, Date and Publisher values are refreshed
#    to reflect the move from Alvearie/ibm.com to LinuxForHealth.
#  - Elements sheet: the stray "Constraint(s)" text that had (incorrectly)
#    been duplicated onto the Extension row is cleared, leaving it only on
#    the Extension.extension row where it belongs.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

$wsMetadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/total-dependent-count"
$wsMetadata.Range("B3").Value = "8.0.0"
$wsMetadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMetadata.Range("B9").Value = "LinuxForHealth Team"

# The Extension.url row's "Fixed Value" column mirrors the canonical URL above.
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/total-dependent-count"

$wsElements.Range("AI2").Value = ""
